$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (1 changes)
$ws.Range("I3").Value = 3.8

# Row 4 (2 changes)
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.62

# Row 5 (5 changes)
$ws.Range("G5").Value = 4.33
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 2.05
$ws.Range("AE5").Value = 23
$ws.Range("AN5").Value = 6

# Row 6 (3 changes)
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7
$ws.Range("AG6").Value = 700

# Row 7 (21 changes)
$ws.Range("H7").Value = 4.9
$ws.Range("L7").Value = 6.2
$ws.Range("N7").Value = 9.75
$ws.Range("P7").Value = 4.9
$ws.Range("Q7").Value = 1.44
$ws.Range("R7").Value = 2.57
$ws.Range("W7").Value = 9.75
$ws.Range("X7").Value = 8
$ws.Range("Z7").Value = 9.5
$ws.Range("AC7").Value = 9.75
$ws.Range("AD7").Value = 10.25
$ws.Range("AE7").Value = 17
$ws.Range("AH7").Value = 25
$ws.Range("AL7").Value = 70
$ws.Range("AM7").Value = 55
$ws.Range("AN7").Value = 3.45
$ws.Range("AO7").Value = 6
$ws.Range("AP7").Value = 13.5
$ws.Range("AU7").Value = 7.6
$ws.Range("AZ7").Value = 250
$ws.Range("BB7").Value = 350

# Row 8 (22 changes)
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 2.92
$ws.Range("K8").Value = 2.07
$ws.Range("L8").Value = 3.5
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 3.7
$ws.Range("Q8").Value = 1.72
$ws.Range("R8").Value = 2.05
$ws.Range("T8").Value = 2.72
$ws.Range("Y8").Value = 8.75
$ws.Range("AB8").Value = 22
$ws.Range("AD8").Value = 6.6
$ws.Range("AH8").Value = 10.75
$ws.Range("AI8").Value = 16.5
$ws.Range("AL8").Value = 23
$ws.Range("AM8").Value = 27
$ws.Range("AT8").Value = 2.72
$ws.Range("AV8").Value = 60
$ws.Range("AX8").Value = 16.5
$ws.Range("AY8").Value = 23
$ws.Range("BA8").Value = 110
$ws.Range("BB8").Value = 300

# Row 9 (33 changes)
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3.55
$ws.Range("I9").Value = 3.35
$ws.Range("J9").Value = 2.57
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 3.8
$ws.Range("S9").Value = 1.36
$ws.Range("T9").Value = 2.92
$ws.Range("U9").Value = 1.6
$ws.Range("V9").Value = 2.2
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 10.75
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 18.5
$ws.Range("AA9").Value = 14.5
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 12.5
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 19.5
$ws.Range("AJ9").Value = 11.5
$ws.Range("AK9").Value = 45
$ws.Range("AL9").Value = 27
$ws.Range("AM9").Value = 30
$ws.Range("AN9").Value = 4.05
$ws.Range("AO9").Value = 10
$ws.Range("AQ9").Value = 37
$ws.Range("AT9").Value = 2.92
$ws.Range("AU9").Value = 6.8
$ws.Range("AW9").Value = 5.4
$ws.Range("AX9").Value = 18
$ws.Range("AY9").Value = 23
$ws.Range("AZ9").Value = 90
$ws.Range("BA9").Value = 110

# Row 11 (42 changes)
$ws.Range("G11").Value = 1.45
$ws.Range("H11").Value = 4.2
$ws.Range("I11").Value = 6.5
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 2.22
$ws.Range("L11").Value = 6.2
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.3
$ws.Range("Q11").Value = 1.87
$ws.Range("R11").Value = 1.87
$ws.Range("S11").Value = 1.4
$ws.Range("T11").Value = 2.72
$ws.Range("U11").Value = 2.05
$ws.Range("V11").Value = 1.7
$ws.Range("W11").Value = 6.2
$ws.Range("X11").Value = 6.3
$ws.Range("Z11").Value = 9.25
$ws.Range("AA11").Value = 12.5
$ws.Range("AC11").Value = 7.5
$ws.Range("AD11").Value = 8.25
$ws.Range("AE11").Value = 21
$ws.Range("AF11").Value = 110
$ws.Range("AH11").Value = 15
$ws.Range("AI11").Value = 37
$ws.Range("AJ11").Value = 21
$ws.Range("AK11").Value = 150
$ws.Range("AL11").Value = 75
$ws.Range("AM11").Value = 80
$ws.Range("AN11").Value = 3.15
$ws.Range("AO11").Value = 6.8
$ws.Range("AP11").Value = 19
$ws.Range("AQ11").Value = 21
$ws.Range("AR11").Value = 60
$ws.Range("AT11").Value = 2.72
$ws.Range("AU11").Value = 8.75
$ws.Range("AW11").Value = 7.6
$ws.Range("AX11").Value = 40
$ws.Range("AY11").Value = 45
$ws.Range("AZ11").Value = 300
$ws.Range("BA11").Value = 350

# Row 12 (14 changes)
$ws.Range("I12").Value = 11.5
$ws.Range("J12").Value = 1.57
$ws.Range("K12").Value = 2.87
$ws.Range("R12").Value = 2.92
$ws.Range("T12").Value = 3.9
$ws.Range("W12").Value = 10.5
$ws.Range("X12").Value = 7.6
$ws.Range("Z12").Value = 7.9
$ws.Range("AI12").Value = 110
$ws.Range("AK12").Value = 400
$ws.Range("AN12").Value = 3.35
$ws.Range("AT12").Value = 3.9
$ws.Range("AW12").Value = 11.75
$ws.Range("AY12").Value = 45
